# "Address lookup test data and locators fix" -
# replace the long, pre-formatted full address in the
# "recipientfulladdress2" test-data column (Q) with a clean street
# address, then tidy up the column width / selection that Excel
# recalculates as a side effect of that edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Q2 held "Pennsylvania Turnpike Commission, 700 S Eisenhower Blvd,
# Middletown, PA 17057" - replace it with the shorter street address.
$ws.Range("Q2").Value = "700 South Eisenhower Boulevard"

# Setting .Value resets any quote-prefix/number formatting on the cell;
# re-apply the formatting (borders etc.) that the rest of the row uses,
# copying it from the neighboring address cell.
$ws.Range("P2").Copy()
$ws.Range("Q2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# The column was sized to fit the old, much longer address - shrink it
# back down now that Q2 holds a shorter string.
$ws.Columns("Q:Q").ColumnWidth = 27.67

# Leave the same cell selected as in the authored workbook.
$ws.Range("C2").Select()
